$wb = $excel.ActiveWorkbook

# --- Sheet "1_1": update timestamp text ---
$ws1 = $wb.Worksheets.Item("1_1")
$ws1.Range("C2").Value = "12/31/2023, 00:23:14"
$ws1.Range("D7").Value = "12/31/2023, 00:23:14"

# --- Sheet "summary": rename products and re-sort Makespan values ---
$ws2 = $wb.Worksheets.Item("summary")
$ws2.Range("A2").Value = "PV-1_PI-1"
$ws2.Range("A3").Value = "PV-2_PI-1"
$ws2.Range("A4").Value = "PV-2_PI-2"
$ws2.Range("A5").Value = "PV-3_PI-1"

$ws2.Range("B2").Value = 5000
$ws2.Range("B5").Value = 1000
